$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Replace the phone-number-like chat ids in A2:A20 with a simple 1..10
#     sequence, dropping rows 12:20 entirely (delete whole rows) and
#     clearing the custom font/style that had been applied to A2:A20.

# Remove the extra formatting (font) that was applied to A2:A20 so the
# cells fall back to the default style.
$ws.Range("A2:A20").ClearFormats()

# Drop rows 12-20 completely (shifts nothing below them, just removes
# them from the sheet - dimension shrinks to A1:A11).
$ws.Range("A12:A20").EntireRow.Delete()

# Renumber the remaining ids as a plain 1-10 sequence.
for ($i = 2; $i -le 11; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 1
}

# Move the selection/active cell to A11 (single cell, not a range).
$ws.Range("A11").Select() | Out-Null

# Drop the explicit printer/paper-size page setup that had been carried
# over from the source document (reset to "no explicit paper size").
$ws.PageSetup.PaperSize = $null
